# Update "想去人数" (F column) values on the "展览" and "全部类型" sheets
# to reflect newly generated counts (gh-pages output regenerated at 456a3b4).

$wb = $excel.ActiveWorkbook

# Map of sheet name -> { row number -> new F value }
$updates = @{
    "展览" = @{
        2  = 286
        3  = 48
        4  = 3506
        5  = 2186
        6  = 427
        9  = 53
        10 = 1297
        12 = 1736
        13 = 129
    }
    "全部类型" = @{
        2  = 286
        3  = 48
        4  = 3506
        5  = 2186
        6  = 427
        10 = 53
        13 = 1297
        15 = 1736
        16 = 129
    }
}

foreach ($sheetName in $updates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $rows = $updates[$sheetName]
    foreach ($row in $rows.Keys) {
        $ws.Range("F$row").Value = $rows[$row]
    }
}
